$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.485.82"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "2.058.32"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.30"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.94"
$ws.Range("E8").Value = "  -4.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.57"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.359"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.895"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.74"
$ws.Range("E14").Value = "  -6.46%  "
$ws.Range("D15").Value = "2.358.83"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("E16").Value = "  -7.31%  "
$ws.Range("D17").Value = "2.059.16"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "36.437.06"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.78"
$ws.Range("E19").Value = "  -9.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.05"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("E21").Value = "  -5.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.51"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -5.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.89"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.20"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -7.63%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  -7.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0593"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  -6.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.88"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0215"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.84"
$ws.Range("E42").Value = "  -8.60%  "
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.75"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0902"
$ws.Range("E45").Value = "  -10.27%  "
$ws.Range("D46").Value = "1.390.49"
$ws.Range("E46").Value = "  +6.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.70"
$ws.Range("E47").Value = "  -8.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.37"
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").Value = "2.247.27"
$ws.Range("E51").Value = "  +0.21%  "
